$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet updates ---
$wsSchedule.Range("B2").Value = 46041.22916666666
$wsSchedule.Range("C2").Value = 5.5
$wsSchedule.Range("D2").Value = 20.79
$wsSchedule.Range("E2").Value = 509.8020524999999
$wsSchedule.Range("F2").Value = 24.52150324675324
$wsSchedule.Range("A4").Value = 46041.95833333334
$wsSchedule.Range("B4").Value = 46042.125
$wsSchedule.Range("E4").Value = 444.8536169999999
$wsSchedule.Range("F4").Value = 29.42153551587301
$wsSchedule.Range("A5").Value = 46042.29166666666
$wsSchedule.Range("C5").Value = 9
$wsSchedule.Range("D5").Value = 34.02
$wsSchedule.Range("E5").Value = -73.51220174999999
$wsSchedule.Range("F5").Value = -2.160852491181658

# --- Detailed sheet updates ---
$wsDetailed.Range("E13").Value = "OFF"
$wsDetailed.Range("B41").Value = 64.02373
$wsDetailed.Range("B42").Value = 73.20007
$wsDetailed.Range("B43").Value = 59.9853
$wsDetailed.Range("B44").Value = 65
$wsDetailed.Range("C44").Value = "historical"
$wsDetailed.Range("B45").Value = 59.14887
$wsDetailed.Range("B46").Value = 59.46714
$wsDetailed.Range("B47").Value = 58.71986
$wsDetailed.Range("B48").Value = 57.08
$wsDetailed.Range("E48").Value = "ON"
$wsDetailed.Range("B49").Value = 57.06003
$wsDetailed.Range("B52").Value = 56.98
$wsDetailed.Range("B53").Value = 56.98
$wsDetailed.Range("B54").Value = 56.98
$wsDetailed.Range("E56").Value = "OFF"
$wsDetailed.Range("B60").Value = 64.89
$wsDetailed.Range("B61").Value = 64.89
$wsDetailed.Range("B62").Value = 73.20005
$wsDetailed.Range("B63").Value = 61.4478
$wsDetailed.Range("B64").Value = 36.06
$wsDetailed.Range("E64").Value = "ON"
$wsDetailed.Range("B65").Value = 0.98766
$wsDetailed.Range("B67").Value = 0.51
$wsDetailed.Range("B68").Value = -1.17721
$wsDetailed.Range("B69").Value = -6
$wsDetailed.Range("B70").Value = -7.12181
$wsDetailed.Range("B71").Value = -7.42877
$wsDetailed.Range("B72").Value = -7.78472
$wsDetailed.Range("B73").Value = -7.73401
$wsDetailed.Range("B74").Value = -7.82489
$wsDetailed.Range("B75").Value = -9.99
$wsDetailed.Range("B77").Value = -12.01
$wsDetailed.Range("B78").Value = -11.01
$wsDetailed.Range("B79").Value = -10.39703
$wsDetailed.Range("B80").Value = -8.38635
$wsDetailed.Range("B81").Value = -6.8
$wsDetailed.Range("B82").Value = -5.51
$wsDetailed.Range("B83").Value = -5.56366
$wsDetailed.Range("B84").Value = -1.01654
$wsDetailed.Range("B85").Value = 9.338240000000001
$wsDetailed.Range("B86").Value = 9.46283
$wsDetailed.Range("B87").Value = 33.1597
$wsDetailed.Range("B88").Value = 57.3
$wsDetailed.Range("B89").Value = 57.3
$wsDetailed.Range("B90").Value = 60.00844
$wsDetailed.Range("B91").Value = 57.3
$wsDetailed.Range("B92").Value = 57.3
$wsDetailed.Range("B93").Value = 57.3
$wsDetailed.Range("B94").Value = 57.06007
$wsDetailed.Range("B95").Value = 57.3
$wsDetailed.Range("B96").Value = 57.3
$wsDetailed.Range("B97").Value = 57.3
